# "Generate Report for Archive"
#
# Two changes in the localization-status report:
#   1. Status text "Ready for handoff" -> "In Translation" (Status column on
#      every sheet: Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4 all share the
#      same cached string).
#   2. The Status column is narrowed (shrinks from ~17.22 chars to ~13.41
#      chars) on all three sheets.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- 1. Update the status value everywhere it appears -----------------

$overview = $wb.Worksheets.Item("Overview")
foreach ($addr in @("E2", "F2", "E3", "F3", "E4", "F4")) {
    $overview.Range($addr).Value = $newStatus
}

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in @("C2", "C3", "C4")) {
        $ws.Range($addr).Value = $newStatus
    }
}

# --- 2. Narrow the Status column on all three sheets -------------------
# Target stored width is ~13.4101845877511 characters; the host quantizes
# ColumnWidth writes to 1/6-character increments, so 12.5 (mid-bucket) is
# the input that lands on the nearest reachable grid value. Note:
# Columns.Item("E") (letter index) isn't supported by this host, so the
# whole-column Range("E:E") form is used instead.

$overview.Range("E:E").ColumnWidth = 12.5
$overview.Range("F:F").ColumnWidth = 12.5

$wb.Worksheets.Item("zh-cn").Range("C:C").ColumnWidth = 12.5
$wb.Worksheets.Item("de-de").Range("C:C").ColumnWidth = 12.5

Write-Output "edit complete"
